# Regenerate merged AHB files
# - Rename the "_old"/"_new" header-suffix columns to "_FV2210"/"_FV2304"
# - Wrap the data range in an Excel Table (Table1) with an AutoFilter
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- 1. Rename header cells -------------------------------------------------
# Columns A-J: "<Name>_old"  -> "<Name>_FV2210"
$headersFV2210 = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
for ($i = 0; $i -lt $headersFV2210.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2210[$i]
}

# Column K ("diff") is unchanged.

# Columns L-U: "<Name>_new" -> "<Name>_FV2304"
$headersFV2304 = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
for ($i = 0; $i -lt $headersFV2304.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2304[$i]
}

# --- 2. Turn the data range into an Excel Table (adds autoFilter too) ------
$tableRange = $ws.Range("A1:U75")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
